# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 185
$wsExhibit.Range("F4").Value = 343
$wsExhibit.Range("F5").Value = 410
$wsExhibit.Range("F6").Value = 258
$wsExhibit.Range("F7").Value = 2389
$wsExhibit.Range("F8").Value = 407
$wsExhibit.Range("F9").Value = 6151
$wsExhibit.Range("F10").Value = 157
$wsExhibit.Range("F11").Value = 394

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 185
$wsAll.Range("F4").Value = 343
$wsAll.Range("F5").Value = 410
$wsAll.Range("F6").Value = 258
$wsAll.Range("F9").Value = 2389
$wsAll.Range("F10").Value = 407
$wsAll.Range("F11").Value = 6151
$wsAll.Range("F12").Value = 157
$wsAll.Range("F13").Value = 394
